$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = 1.88
$ws.Range("U2").Value = 2.02
$ws.Range("AG2").Value = 9.6

$ws.Range("F5").Value = 5.9
$ws.Range("G5").Value = 7.2
$ws.Range("H5").Value = 1.54
$ws.Range("I5").Value = 1.64
$ws.Range("W5").Value = 1.16
$ws.Range("Z5").Value = 10.5
$ws.Range("AA5").Value = 1000
$ws.Range("AG5").Value = 26
$ws.Range("AH5").Value = 23
$ws.Range("AK5").Value = 90
$ws.Range("AO5").Value = 9.800000000000001

$ws.Range("G6").Value = 1.46
$ws.Range("J6").Value = 4.7

$ws.Range("F7").Value = 1.73
$ws.Range("G7").Value = 2.02
$ws.Range("P7").Value = 1.85
$ws.Range("W7").Value = 1.98

$ws.Range("F8").Value = 1.42
$ws.Range("G8").Value = 1.43
$ws.Range("J8").Value = 5.1
$ws.Range("K8").Value = 5.2
$ws.Range("V8").Value = 1.11
$ws.Range("W8").Value = 3.35
$ws.Range("AN8").Value = 7

$ws.Range("H10").Value = 1.97
$ws.Range("K10").Value = 5.4
$ws.Range("N10").Value = 2.68
$ws.Range("Q10").Value = 1.99
$ws.Range("W10").Value = 1.33

$ws.Range("F11").Value = 6.8
$ws.Range("I11").Value = 1.62
$ws.Range("O11").Value = 1.39
$ws.Range("U11").Value = 1.78
$ws.Range("V11").Value = 2.6
$ws.Range("AN11").Value = 190

$ws.Range("F12").Value = 2.86
$ws.Range("G12").Value = 2.9
$ws.Range("H12").Value = 2.84
$ws.Range("I12").Value = 2.86
$ws.Range("P12").Value = 1.7
$ws.Range("Q12").Value = 2.36
$ws.Range("U12").Value = 1.97
$ws.Range("Y12").Value = 9.4
$ws.Range("AC12").Value = 7.2

$ws.Range("G13").Value = 4.7
$ws.Range("I13").Value = 2.08
$ws.Range("J13").Value = 3.35
$ws.Range("W13").Value = 1.28
$ws.Range("AI13").Value = 980

$ws.Range("G14").Value = 3.7
$ws.Range("T14").Value = 2
$ws.Range("AA14").Value = 980
$ws.Range("AE14").Value = 980
$ws.Range("AO14").Value = 980

$ws.Range("F16").Value = 2.28
$ws.Range("G16").Value = 2.3
$ws.Range("H16").Value = 3.6
$ws.Range("I16").Value = 3.65
$ws.Range("L16").Value = 1.41
$ws.Range("V16").Value = 1.37
$ws.Range("W16").Value = 1.76
$ws.Range("Y16").Value = 13.5
$ws.Range("Z16").Value = 24
$ws.Range("AB16").Value = 9.6
$ws.Range("AD16").Value = 14.5

$ws.Range("F17").Value = 1.78
$ws.Range("I17").Value = 5.8
$ws.Range("J17").Value = 3.8
$ws.Range("K17").Value = 3.9
$ws.Range("R17").Value = 1.33
$ws.Range("T17").Value = 1.97
$ws.Range("AA17").Value = 150

$ws.Range("J18").Value = 3.35
$ws.Range("M18").Value = 1.09
$ws.Range("N18").Value = 3.2
$ws.Range("O18").Value = 1.38
$ws.Range("R18").Value = 1.29
$ws.Range("S18").Value = 3.95
$ws.Range("T18").Value = 1.04
$ws.Range("U18").Value = 1.04
$ws.Range("X18").Value = 12.5
$ws.Range("Y18").Value = 8.800000000000001
$ws.Range("Z18").Value = 980
$ws.Range("AA18").Value = 980
$ws.Range("AB18").Value = 13
$ws.Range("AC18").Value = 7.8
$ws.Range("AD18").Value = 11
$ws.Range("AE18").Value = 980
$ws.Range("AF18").Value = 980
$ws.Range("AG18").Value = 16.5
$ws.Range("AH18").Value = 990
$ws.Range("AI18").Value = 980
$ws.Range("AJ18").Value = 80
$ws.Range("AK18").Value = 55
$ws.Range("AL18").Value = 65
$ws.Range("AM18").Value = 150
$ws.Range("AO18").Value = 22

$ws.Range("F19").Value = 3.95
$ws.Range("G19").Value = 4.1
$ws.Range("H19").Value = 2.08
$ws.Range("I19").Value = 2.1
$ws.Range("L19").Value = 1.34
$ws.Range("P19").Value = 2.16
$ws.Range("V19").Value = 1.9
$ws.Range("W19").Value = 1.33
$ws.Range("X19").Value = 16.5
$ws.Range("Y19").Value = 11
$ws.Range("Z19").Value = 13.5
$ws.Range("AB19").Value = 16.5
$ws.Range("AD19").Value = 10.5
$ws.Range("AE19").Value = 20
$ws.Range("AF19").Value = 29
$ws.Range("AG19").Value = 15.5
$ws.Range("AI19").Value = 32
$ws.Range("AK19").Value = 44
$ws.Range("AL19").Value = 48
$ws.Range("AN19").Value = 38

$ws.Range("L20").Value = 1.22
$ws.Range("M20").Value = 1.02
$ws.Range("N20").Value = 6
$ws.Range("O20").Value = 1.12
$ws.Range("Q20").Value = 1.41
$ws.Range("R20").Value = 1.79
$ws.Range("S20").Value = 1.89
$ws.Range("T20").Value = 1.01
$ws.Range("U20").Value = 1.01
$ws.Range("V20").Value = 1.04
$ws.Range("W20").Value = 6.8
$ws.Range("X20").Value = 46
$ws.Range("Y20").Value = 85
$ws.Range("Z20").Value = 1000
$ws.Range("AA20").Value = 1000
$ws.Range("AB20").Value = 17
$ws.Range("AC20").Value = 30
$ws.Range("AD20").Value = 100
$ws.Range("AE20").Value = 1000
$ws.Range("AF20").Value = 10.5
$ws.Range("AG20").Value = 15.5
$ws.Range("AH20").Value = 55
$ws.Range("AI20").Value = 1000
$ws.Range("AJ20").Value = 10
$ws.Range("AK20").Value = 18
$ws.Range("AL20").Value = 55
$ws.Range("AM20").Value = 1000
$ws.Range("AN20").Value = 1000
$ws.Range("AO20").Value = 1000

